# Update column G ("K") values on Sheet1 per regenerated save_data
# (commit: regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$gValues = @{
    2  = 3
    3  = 0
    4  = 2
    5  = 3
    6  = 1
    7  = 0
    8  = 2
    9  = 1
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 0
    15 = 3
    16 = 1
    17 = 2
    18 = 2
    19 = 1
    20 = 0
    21 = 5
    22 = 0
    23 = 3
    24 = 0
    25 = 2
    26 = 1
    27 = 0
    28 = 1
    29 = 1
    30 = 2
    31 = 1
    32 = 0
    33 = 2
    34 = 1
    35 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
